$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change "Helado" to "helado" in C10
$ws.Range("C10").Value = "helado"

# Clear the product row that was removed (row 11: B11:E11)
$ws.Range("B11:E11").Value = ""
